$d = $word.ActiveDocument

function Replace-Text($oldText, $newText) {
    $d.Content.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}

# Phase 1: move each original value to a unique placeholder token so that the
# subsequent phase-2 replacements cannot accidentally match text that was
# just inserted by an earlier step (the 8 values form a rotation, so naive
# sequential replacement could cause collisions).
Replace-Text "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas." "@@SLOT_1@@"
Replace-Text "5840938 - Marcelo Rodrigues de Holanda" "@@SLOT_2@@"
Replace-Text "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso." "@@SLOT_3@@"
Replace-Text "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental." "@@SLOT_4@@"
Replace-Text "Aula expositiva e exercícios dirigidos." "@@SLOT_5@@"
Replace-Text "Média ponderada de exercícios e provas." "@@SLOT_6@@"
Replace-Text "Prova única com nota igual ou superior a 5,0." "@@SLOT_7@@"
Replace-Text "Estudos de caso: EPIA de origem." "@@SLOT_8@@"

# Phase 2: replace each placeholder with the value that now belongs there,
# per the diff (value from the "next" slot in the rotation).
Replace-Text "@@SLOT_1@@" "Estudos de caso: avaliar a singularidade do caso apresentado, com a necessária identificação das etapas e das peculiaridades que os fizeram próprios a serem aplicados em estudo de caso."
Replace-Text "@@SLOT_2@@" "Avaliar casos de impacto ambiental que proporcionaram singularidades em algumas de suas etapas."
Replace-Text "@@SLOT_3@@" "Estudos de casos específicos e as etapas necessárias na avaliação de um impacto ambiental."
Replace-Text "@@SLOT_4@@" "Aula expositiva e exercícios dirigidos."
Replace-Text "@@SLOT_5@@" "Média ponderada de exercícios e provas."
Replace-Text "@@SLOT_6@@" "Prova única com nota igual ou superior a 5,0."
Replace-Text "@@SLOT_7@@" "Estudos de caso: EPIA de origem."
Replace-Text "@@SLOT_8@@" "5840938 - Marcelo Rodrigues de Holanda"
